# Historias de Usuario - add "Registrar Personal" user story row.
#
# Current layout (before):
#   Row3: Administrador del sistema | Registrar Usuario  | Tener los usuarios en la base de datos      | ALTA
#   Row4: Usuario Admin/Invitado    | Iniciar Sesion      | Identificar el rol...                       | ALTA
#   Row5: Invitado                  | Control de Entrada/Salida | Tener registro de las personas...     | ALTA
#   Row6: Administrador del sistema | Reportes            | Dar informacion al Administrador...          | ALTA
#   Row7: Usuario Admin/Invitado    | Diseño              | Dar una buena imagen del aplicativo...       | MEDIA
#
# Target layout (after): a new row is inserted after row 3 ("Registrar
# Usuario") for a new "Registrar Personal" story, pushing the rest down by
# one row, and the "Registrar Usuario" need text is reworded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reword the existing "Registrar Usuario" need (row 3) to be more specific.
$ws.Range("D3").Value = "Tener los usuarios en la base de datos para el ingreso del aplicativo"

# Insert a new row at position 4 - this shifts the old rows 4-7 down to 5-8,
# carrying their formatting (including row heights) along automatically.
$ws.Rows.Item(4).Insert()

# Clone row 3's cell formatting into the freshly inserted row 4 so the new
# row matches the look of the rest of the table (borders/fill/font/etc.).
$ws.Range("B3:F3").Copy()
$ws.Range("B4:F4").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4104)
$ws.Application.CutCopyMode = $false

# Fill in the new "Registrar Personal" story.
$ws.Range("B4").Value = "Administrador del sistema"
$ws.Range("C4").Value = "Registrar Personal"
$ws.Range("D4").Value = "Tener el personal en la base de datos para el control de entrada y salida"
$ws.Range("E4").Value = "ALTA"

# Row 3 keeps its taller "wrapped need" height; row 4 (new) matches it.
$ws.Rows.Item(3).RowHeight = 51
$ws.Rows.Item(4).RowHeight = 51

# Update the selection to reflect where the author left off editing.
$ws.Range("D5").Select()
